# 27 Feb-Added Script for external user profile
#
# - RunToTest flag for row 2 (Valid_Login) changed from Y to N.
# - New row 5 added with a new external user profile:
#     RunToTest = Y, UserName = Nasreenk@winjit.com, Password = Nasreenk@winjit.com
#   Both UserName and Password cells get a mailto hyperlink + Hyperlink style,
#   matching the existing rows in the sheet.
# - Selection moved to B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: RunToTest Y -> N
$ws.Range("A2").Value = "N"

# New row 5: external user profile
$ws.Range("A5").Value = "Y"
$ws.Range("B5").Value = "Nasreenk@winjit.com"
$ws.Range("C5").Value = "Nasreenk@winjit.com"

$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:Nasreenk@winjit.com")
$ws.Range("B5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Nasreenk@winjit.com")
$ws.Range("C5").Style = "Hyperlink"

# Move the sheet selection to B9
$ws.Range("B9").Select()
